$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.896.53"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.862.03"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.67"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5044"
$ws.Range("E7").Value = "  -1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3635"
$ws.Range("E8").Value = "  -2.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07160"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8929"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.73"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").Value = "1.857.58"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.44"
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.224"
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008483"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.17"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9996"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "26.933.76"
$ws.Range("E20").Value = "  -1.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.027"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "2.081.99"
$ws.Range("E22").Value = "  -1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.37"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.394"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.07"
$ws.Range("E25").Value = "  -1.94%  "
$ws.Range("E26").Value = "  -3.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.86"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.075"
$ws.Range("E28").Value = "  -1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.04"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.668"
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09246"
$ws.Range("E32").Value = "  +2.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05092"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7506"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.997"
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.150"
$ws.Range("E36").Value = "  -0.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.264"
$ws.Range("E37").Value = "  +6.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.526"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01996"
$ws.Range("E39").Value = "  -2.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5576"
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.069"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "118.57"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.530"
$ws.Range("E43").Value = "  -0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.517"
$ws.Range("E44").Value = "  +2.19%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1468"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4686"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.0000"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.05"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.561"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.74"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.87"
$ws.Range("E51").Value = "  -2.34%  "
